# Update batch_import template: reorder KEY/GROUPS columns, add COMPLIANCE_GROUPS
# column, and fill in sample values for the new "anssi" compliance group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sshKey = @"
---- BEGIN SSH2 PUBLIC KEY ----
Comment: example
AAAAB3NzaC1kc3MAAACBAMXXH+SzAIPRN38GehSARboF873Ic5utBjMcXx1IuFNTvvBi2j
fCyjCBqb66sgS8pdzUl+YyH4sMDp67Q9RKI9po3ePtV03rldPQjtqgmXt2B2eQ6SKXDO3g
+nN9LLEVXp9MpS7g9VnyDWUQCAxu+Khp+uZDzwSy7IVxRm/HHU2dAAAAFQCl1FWQ7bTyoY
7RtEvB6rhqGyY/8QAAAIEAxYgBAFfVKvSC3AZkwWuB4hPLlBeKhL4Yt87vblimHWlaOSFU
llKnCGmdc7R2NL3JZFP210yjapZY25YTKpkO8pdavazVqbzBd1EEtZ93umDqWua2yqPOc8
6MoZJbk7OTJjZRlpd1XZwSI3XgyxaDtf+tCh14ikG13k4A1iKd3/MAAACBALHyHX29XFe3
VseZeG+CiYMfc3qXbMQgpWdZeopg/1Z3qw46Kx4iiNgtZcB7BdoYdIhDvTu+xkffbG22h9
YQnxyM9Kz/cqjKdKHp2VBX/IJU4vEkIPF+kdFPToLvJc+qkIvd1kDqUUW+e6dD6PkpNDdh
gOn/vcgro4IwufBActyG
---- END SSH2 PUBLIC KEY ----
"@

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "HOST"
$ws.Range("B1").Value = "PORT"
$ws.Range("C1").Value = "TYPE"
$ws.Range("D1").Value = "NODE"
$ws.Range("E1").Value = "KEY"
$ws.Range("F1").Value = "USERNAME"
$ws.Range("G1").Value = "PASSWORD"
$ws.Range("H1").Value = "GROUPS"
$ws.Range("I1").Value = "COMPLIANCE_GROUPS"

# --- Row 2 --------------------------------------------------------------
$ws.Range("A2").Clear()
$ws.Range("B2").Value = 22
$ws.Range("C2").Value = "CbwRam::RemoteAccess::Ssh::WithPassword"
$ws.Range("D2").Value = "master"
$ws.Range("E2").Clear()
$ws.Range("F2").Value = "admin"
$ws.Range("G2").Value = "SuperPassword"
$ws.Range("H2").Value = "production, test"
$ws.Range("I2").Value = "anssi"

# --- Row 3 --------------------------------------------------------------
$ws.Range("A3").Value = "server02.example.com"
$ws.Range("B3").Value = 22
$ws.Range("C3").Value = "CbwRam::RemoteAccess::Ssh::WithPassword"
$ws.Range("D3").Value = "master"
$ws.Range("E3").Value = $sshKey
$ws.Range("E3").WrapText = $true
$ws.Range("F3").Value = "user"
$ws.Range("G3").Clear()
$ws.Range("H3").Clear()
$ws.Range("I3").Value = "anssi"

# --- Row 4 --------------------------------------------------------------
$ws.Range("A4").Value = "127.0.1.1"
$ws.Range("B4").Value = 5985
$ws.Range("C4").Value = "CbwRam::RemoteAccess::WinRm::WithNegotiate"
$ws.Range("D4").Value = "slave1"
$ws.Range("E4").Clear()
$ws.Range("F4").Value = "Administrator"
$ws.Range("G4").Value = "Nu6K2WBX"
$ws.Range("H4").Value = "preproduction"
$ws.Range("I4").Value = "anssi"

# --- Row heights ----------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(2).RowHeight = 13.8
$ws.Rows.Item(3).RowHeight = 163.5
$ws.Rows.Item(4).RowHeight = 13.8

# --- Column widths ----------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 15.75
$ws.Columns.Item(5).ColumnWidth = 6.252
$ws.Columns.Item(6).ColumnWidth = 11.752
$ws.Columns.Item(7).ColumnWidth = 14.42
$ws.Columns.Item(8).ColumnWidth = 13.584
$ws.Columns.Item(9).ColumnWidth = 19.084

# --- Selection / view -------------------------------------------------------
$ws.Range("D6").Select() | Out-Null
